# Generate Report for Handoff
# The e8d00c08-...md file has finished translation/review and is now ready
# to be handed off: update its Status/Priority/Handoff-timestamp on every
# sheet that tracks it (the per-locale sheets + the Overview roll-up).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 corresponds to e8d00c08-...md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-25 02:14:09"

# --- de-de sheet: row 3 corresponds to e8d00c08-...md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-25 02:14:14"

# --- Overview sheet: row 3 corresponds to e8d00c08-...md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 02:14:14"
